# Weekly update: insert the newest week's row for Brócoli at
# Macroferia Regional de Talca, pushing all prior rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 454:476 down to 455:477, inserting a fresh blank row at 454.
$ws.Rows("454:454").Insert()

# Populate the newly inserted row 454 with this week's values.
$ws.Range("A454").Value = 5
$ws.Range("B454").Value = "Macroferia Regional de Talca"
$ws.Range("C454").Value = "Maule"
$ws.Range("D454").Value = 44939
$ws.Range("E454").Value = 7
$ws.Range("F454").Value = 100112023
$ws.Range("G454").Value = "Brócoli"
$ws.Range("H454").Value = "Sin especificar"
$ws.Range("I454").Value = "Primera"
$ws.Range("J454").Value = 4000
$ws.Range("K454").Value = 800
$ws.Range("L454").Value = 800
$ws.Range("M454").Value = 800
$ws.Range("N454").Value = '$/unidad'
$ws.Range("O454").Value = "Región del Maule"
$ws.Range("P454").Value = 800
$ws.Range("Q454").Value = 1
$ws.Range("R454").Value = "Hortaliza"
